# Login/password section restructuring:
#  - "Administrador" paragraph loses its trailing _GoBack bookmark.
#  - A new paragraph with leading spaces + "Usuario: 1M5533yuc" is inserted.
#  - A new paragraph "Contra: prueba" is inserted, now carrying the _GoBack
#    bookmark and an underline on its paragraph mark run.
#  - The former "Usuario:/Contra:" paragraphs are replaced by a single
#    "Selección : Adminstrador" paragraph (with proofErr spell/gram markers),
#    reusing the original paragraph's identity.
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le ($count - 2); $i++) {
    $t1 = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($t1 -ne "Administrador") { continue }
    $t2 = $d.Paragraphs($i + 1).Range.Text.TrimEnd([char]13)
    $t3 = $d.Paragraphs($i + 2).Range.Text.TrimEnd([char]13)
    if ($t2 -eq "Usuario: 1M5533yuc" -and $t3 -eq "Contra: prueba") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the Administrador/Usuario/Contra paragraph block"
}

$pAdmin = $d.Paragraphs($targetIndex)
$pContra = $d.Paragraphs($targetIndex + 2)
$target = $d.Range($pAdmin.Range.Start, $pContra.Range.End)

$pkg = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C40B7E" w:rsidRPr="00C40B7E" w:rsidRDefault="00C40B7E" w:rsidP="00074DBE"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Administrador</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">         </w:t></w:r><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Usuario: 1M5533yuc</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:firstLine="708"/><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:u w:val="single"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Contra: prueba</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00090651" w:rsidRDefault="00090651" w:rsidP="00C40B7E"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:firstLine="708"/><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Selección :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="apple-converted-space"/><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="313534"/><w:sz w:val="28"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Adminstrador</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($pkg)
